$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 226, shifting rows 226:327 down to 227:328
$ws.Rows(226).Insert()

# The new row 226 is a duplicate of the (original, still-unedited) row 225
$ws.Range("A225:R225").Copy()
$ws.Range("A226:R226").PasteSpecial()

# Now edit the original row 225 in place: new date and new volume
$ws.Range("D225").Value2 = 44510
$ws.Range("J225").Value2 = 180
